$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (diff-driven). D-column numeric-looking values
# are written with a leading apostrophe so Excel stores them as text (matching
# the source workbook's inlineStr cells) instead of auto-converting to numbers.

$ws.Range("D2").Value = "'62.180.49"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "'2.445.88"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'583.04"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'143.08"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.531"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "'2.441.06"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "'26.49"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "'2.871.87"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'62.002.42"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "'2.434.82"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "'10.78"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'7.19"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "'326.76"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -5.21%  "
$ws.Range("D25").Value = "'65.70"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").Value = "'9.09"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "'600.95"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").Value = "'0.0₃0966"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "'2.568.01"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").Value = "'1.01"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "'1.90"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "'0.135"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "'0.376"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "'152.88"
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("D40").Value = "'18.44"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").Value = "'43.14"
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("E48").Value = "  +19.74%  "
$ws.Range("D49").Value = "'0.601"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'19.79"
$ws.Range("E51").Value = "  -0.10%  "
